# Update column G ("K") values for rows 2-41 on the active worksheet.
# These new values represent a recalculated "K" statistic (replacing the
# previous "Strike#"-derived figures) after regenerating save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 3
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 3
    23 = 0
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 1
    40 = 1
    41 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
